$wb = $excel.ActiveWorkbook

# Sheet "汽車" (Car) - registration date: remove stray space
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("E2").Value = "98年03月24日"

# Sheet "債務" (Debt) - creditor name/address: remove stray space
$wsDebt = $wb.Worksheets.Item("債務")
$wsDebt.Range("D2").Value = "合作金庫商業銀行臺南市北區西門路"

# Sheet "債務" (Debt) - acquisition date: remove stray space
$wsDebt.Range("F2").Value = "89年03月29日"
